$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "Sprint Backlog" header cells to include the sprint length.
$ws.Range("G4").Value = "Sprint Backlog #1 - 2 Weeks"
$ws.Range("G10").Value = "Sprint Backlog #2 - 2 Weeks"

# Match the recorded selection in the saved file.
$ws.Range("J18").Select()
